$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value that was bumped by one day
# (serial 45826 -> 45827, i.e. 2025-06-18 -> 2025-06-19) for every
# existing data row (rows 2 through 43).
for ($row = 2; $row -le 43; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45826) {
        $cell.Value2 = 45827
    }
}
